$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header values change ---
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# --- Row 2 values change ---
$ws.Range("B2").Value = 0.48523255723582637
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = 1.2966861758135664

# --- Row 3 values change ---
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 2.1896109654431584
$ws.Range("D3").Value = 0.90151995643308558
$ws.Range("E3").ClearContents()

# --- Update selection to match new reduced range ---
$ws.Range("B1:E3").Select()
